$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 245-247: date (col A, formatted like existing date cells),
# "nuovi pos." (col B), rolling 7-day sum (col C), and the per-100k rate (col D).
$rows = @(
    @{ Row = 245; Date = 44319; B = 1; C = 4; D = 134.8617666891436 },
    @{ Row = 246; Date = 44320; B = 0; C = 4; D = 134.8617666891436 },
    @{ Row = 247; Date = 44321; B = 0; C = 4; D = 134.8617666891436 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $prevRow = $rowNum - 1

    # Copy the previous row's A-cell (date) formatting so the new cell picks
    # up the same style (border/alignment/date number format) without
    # minting a new style entry, then overwrite with the correct value.
    $ws.Range("A$prevRow").Copy($ws.Range("A$rowNum"))
    $ws.Range("A$rowNum").Value = $r.Date

    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
}

Write-Output "done"
